$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting so
# numeric-looking strings (e.g. "210.87") are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.662.46"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.593.62"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "210.87"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "1.817.91"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "1.586.51"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").Value = "64.43"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "26.646.16"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "207.23"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "6.77"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("E23").Value = "  -2.10%  "
$ws.Range("D24").Value = "8.83"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").Value = "145.57"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "7.17"
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").Value = "0.114"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").Value = "0.665"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "1.278.87"
$ws.Range("E35").Value = "  -3.67%  "
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").Value = "0.837"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "5.40"
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "63.40"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "1.731.06"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "0.908"
$ws.Range("E46").Value = "  +9.35%  "
$ws.Range("D47").Value = "89.93"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.45"
$ws.Range("E51").Value = "  -0.41%  "
